$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.695.76"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "2.301.38"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'301.16"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("D6").Value = "'95.52"
$ws.Range("E6").Value = "  -2.07%  "
$ws.Range("D7").Value = "'0.503"
$ws.Range("E7").Value = "  -1.60%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.495"
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").Value = "'34.65"
$ws.Range("E10").Value = "  -3.17%  "
$ws.Range("D11").Value = "'19.26"
$ws.Range("E11").Value = "  +5.72%  "
$ws.Range("D12").Value = "'0.0787"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "'0.118"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "2.664.54"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "2.308.16"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D18").Value = "42.672.71"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "'12.32"
$ws.Range("E19").Value = "  -3.01%  "
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").Value = "'67.43"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").Value = "'234.91"
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("E24").Value = "  +3.09%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("D27").Value = "'24.58"
$ws.Range("E27").Value = "  -3.57%  "
$ws.Range("D28").Value = "'2.19"
$ws.Range("E28").Value = "  +6.32%  "
$ws.Range("D29").Value = "'163.82"
$ws.Range("E29").Value = "  -1.63%  "
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("D31").Value = "'32.28"
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  -1.38%  "
$ws.Range("D34").Value = "'17.57"
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("D35").Value = "'4.46"
$ws.Range("E35").Value = "  -7.56%  "
$ws.Range("E36").Value = "  -2.25%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("D40").Value = "'2.73"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("D42").Value = "1.966.85"
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'18.98"
$ws.Range("E43").Value = "  +5.19%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'10.51"
$ws.Range("E44").Value = "  +4.89%  "
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("D49").Value = "2.530.16"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").Value = "'52.97"
$ws.Range("E50").Value = "  -2.33%  "
$ws.Range("D51").Value = "'72.10"
$ws.Range("E51").Value = "  +0.04%  "
